# Concrete beam shear check ACI 318-19 -- update flexure/shear results
# and highlight the failing Demand Capacity Ratio row.

$d = $word.ActiveDocument

# --- Table 2: Geometry -------------------------------------------------
$d.Tables.Item(2).Cell(2, 3).Range.Text = "60.0"   # Section height h
$d.Tables.Item(2).Cell(5, 3).Range.Text = "6.44"   # Longitudinal tension rebar As

# --- Table 3: Design forces ---------------------------------------------
$d.Tables.Item(3).Cell(3, 3).Range.Text = "50"     # Shear Vu

# --- Table 4: Limit checks ----------------------------------------------
$d.Tables.Item(4).Cell(2, 3).Range.Text = "0.0"    # Stirrup spacing along length - Value
$d.Tables.Item(4).Cell(2, 5).Range.Text = "27.16"  # Stirrup spacing along length - Max.
$d.Tables.Item(4).Cell(2, 6).Range.Text = "✔️"      # Stirrup spacing along length - Ok?

$d.Tables.Item(4).Cell(3, 3).Range.Text = "0.0"    # Stirrup spacing along width - Value
$d.Tables.Item(4).Cell(3, 5).Range.Text = "54.32"  # Stirrup spacing along width - Max.

$d.Tables.Item(4).Cell(4, 3).Range.Text = "0.0"    # Minimum shear reinforcement - Value
$d.Tables.Item(4).Cell(4, 6).Range.Text = "❌"      # Minimum shear reinforcement - Ok?

$d.Tables.Item(4).Cell(5, 3).Range.Text = "0.0"    # Minimum rebar diameter - Value
$d.Tables.Item(4).Cell(5, 4).Range.Text = "0"      # Minimum rebar diameter - Min.
$d.Tables.Item(4).Cell(5, 6).Range.Text = "✔️"      # Minimum rebar diameter - Ok?

# --- Table 5: Shear reinforcement strength -------------------------------
$d.Tables.Item(5).Cell(2, 3).Range.Text = "0.0"                   # Number of stirrups ns
$d.Tables.Item(5).Cell(3, 3).Range.Text = "0.0"                   # Stirrup diameter db
$d.Tables.Item(5).Cell(4, 3).Range.Text = "0.0"                   # Stirrup spacing s
$d.Tables.Item(5).Cell(5, 3).Range.Text = "54.317142857142855"    # Effective height d
$d.Tables.Item(5).Cell(7, 3).Range.Text = "1.67"                  # Required shear reinforcing Av,req
$d.Tables.Item(5).Cell(8, 3).Range.Text = "0.0"                   # Defined shear reinforcing Av
$d.Tables.Item(5).Cell(9, 3).Range.Text = "0.0"                   # Shear rebar strength ØVs

# --- Table 6: Shear strength ---------------------------------------------
$d.Tables.Item(6).Cell(2, 3).Range.Text = "1086.34"  # Effective shear area Acv
$d.Tables.Item(6).Cell(3, 3).Range.Text = "0.00593"  # Longitudinal reinforcement ratio ρw
$d.Tables.Item(6).Cell(4, 3).Range.Text = "0.794"    # Size modification factor λs
$d.Tables.Item(6).Cell(6, 3).Range.Text = "0.47"     # Concrete effective shear stress kc
$d.Tables.Item(6).Cell(7, 3).Range.Text = "38.64"    # Concrete strength ØVc
$d.Tables.Item(6).Cell(8, 3).Range.Text = "307.51"   # Maximum shear strength ØVmax
$d.Tables.Item(6).Cell(9, 3).Range.Text = "38.64"    # Total shear strength ØVn

# Demand Capacity Ratio row: update value and flag it red (now failing).
$d.Tables.Item(6).Cell(11, 3).Range.Text = "1.29"    # DCR value

$dcrFill = 13551615   # BGR for FFC7CE
$dcrFont = 393372     # BGR for 9C0006

for ($col = 1; $col -le 4; $col++) {
    $cell = $d.Tables.Item(6).Cell(11, $col)
    $cell.Shading.BackgroundPatternColor = $dcrFill
    $r = $cell.Range
    $textRange = $d.Range($r.Start, $r.End - 1)
    $textRange.Font.Color = $dcrFont
}
